# Weekly update of Betarraga price data: a new week's record is inserted
# at row 224 and every subsequent record (previously rows 224-248) shifts
# down by one row, with the former row 248 becoming the new row 249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated values for columns D, I, J, K, L, M, O, P across rows 224-249.
# Row 224 receives the brand-new weekly record; rows 225-249 receive the
# data that used to live one row above (224-248).
$rows = @(
    @(224, 44449, "Primera", 160, 750, 800, 775, "Región del Maule",       155),
    @(225, 44161, "Primera", 270, 700, 750, 719, "Región del Maule",       144),
    @(226, 44438, "Primera", 160, 700, 750, 725, "Provincia de Diguillín", 145),
    @(227, 44438, "Segunda", 300, 600, 650, 625, "Provincia de Diguillín", 125),
    @(228, 44251, "Primera", 120, 600, 650, 625, "Provincia de Diguillín", 125),
    @(229, 44442, "Primera", 160, 750, 800, 775, "Región del Maule",       155),
    @(230, 44435, "Primera", 460, 700, 800, 758, "Provincia de Diguillín", 152),
    @(231, 44435, "Primera", 900, 700, 800, 758, "Región del Maule",       152),
    @(232, 44319, "Primera", 120, 600, 650, 625, "Provincia de Diguillín", 125),
    @(233, 44319, "Segunda",  80, 500, 500, 500, "Provincia de Diguillín", 100),
    @(234, 44175, "Primera", 160, 600, 650, 625, "Región del Maule",       125),
    @(235, 44175, "Segunda", 120, 500, 550, 525, "Región del Maule",       105),
    @(236, 44376, "Primera", 120, 600, 650, 625, "Provincia de Diguillín", 125),
    @(237, 44376, "Segunda", 120, 500, 550, 525, "Provincia de Diguillín", 105),
    @(238, 44279, "Primera", 120, 600, 650, 625, "Provincia de Diguillín", 125),
    @(239, 44279, "Segunda", 120, 500, 550, 525, "Provincia de Diguillín", 105),
    @(240, 44412, "Primera", 160, 600, 650, 625, "Provincia de Diguillín", 125),
    @(241, 44412, "Segunda",  80, 500, 550, 525, "Provincia de Diguillín", 105),
    @(242, 44223, "Primera", 290, 700, 800, 748, "Región del Maule",       150),
    @(243, 44314, "Primera", 120, 600, 650, 625, "Provincia de Diguillín", 125),
    @(244, 44448, "Primera", 300, 750, 800, 775, "Región del Maule",       155),
    @(245, 44167, "Primera", 120, 550, 600, 575, "Provincia de Diguillín", 115),
    @(246, 44238, "Primera", 300, 600, 650, 625, "Provincia de Diguillín", 125),
    @(247, 44238, "Segunda",  80, 500, 550, 525, "Provincia de Diguillín", 105),
    @(248, 44399, "Primera", 300, 650, 700, 675, "Provincia de Diguillín", 135),
    @(249, 44400, "Primera", 300, 700, 750, 725, "Provincia de Diguillín", 145)
)

# Row 249 is a brand-new row; seed the columns that stay constant across
# this whole block of records (A,B,C,E,F,G,H,N,Q,R) with the same values
# already used throughout rows 224-248, then give the Fecha cell the same
# date number format used by the other records in the column.
$ws.Cells.Item(249, 1).Value  = $ws.Cells.Item(248, 1).Value()    # A Mercado (codigo)
$ws.Cells.Item(249, 2).Value  = $ws.Cells.Item(248, 2).Value()    # B Mercado
$ws.Cells.Item(249, 3).Value  = $ws.Cells.Item(248, 3).Value()    # C Region
$ws.Cells.Item(249, 5).Value  = $ws.Cells.Item(248, 5).Value()    # E Codreg
$ws.Cells.Item(249, 6).Value  = $ws.Cells.Item(248, 6).Value()    # F Categoria ID
$ws.Cells.Item(249, 7).Value  = $ws.Cells.Item(248, 7).Value()    # G Categoria
$ws.Cells.Item(249, 8).Value  = $ws.Cells.Item(248, 8).Value()    # H Variedad
$ws.Cells.Item(249, 14).Value = $ws.Cells.Item(248, 14).Value()   # N Unidad de comercializacion
$ws.Cells.Item(249, 17).Value = $ws.Cells.Item(248, 17).Value()   # Q Kg o Unidades
$ws.Cells.Item(249, 18).Value = $ws.Cells.Item(248, 18).Value()   # R Clasificacion
$ws.Cells.Item(249, 4).NumberFormat = $ws.Cells.Item(248, 4).NumberFormat()

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 4).Value  = $r[1]   # D Fecha
    $ws.Cells.Item($rowNum, 9).Value  = $r[2]   # I Calidad
    $ws.Cells.Item($rowNum, 10).Value = $r[3]   # J Volumen
    $ws.Cells.Item($rowNum, 11).Value = $r[4]   # K Precio minimo
    $ws.Cells.Item($rowNum, 12).Value = $r[5]   # L Precio maximo
    $ws.Cells.Item($rowNum, 13).Value = $r[6]   # M Precio promedio ponderado
    $ws.Cells.Item($rowNum, 15).Value = $r[7]   # O Origen
    $ws.Cells.Item($rowNum, 16).Value = $r[8]   # P Precio $/Kg
}
